# Renamed rain model constants
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A19").Value = "fi_lidar_rain_reflectivity"
$ws.Range("A20").Value = "fi_lidar_rain_intensity"

$ws.Range("A21").Select()
